# error solve ifrs list
# Corrects the IFRS-consolidated financial figures in the "company_list"
# sheet (rows 2-9, i.e. fiscal years 2014-2021E) which had been populated
# with wrong (shifted/duplicated) numbers. Also removes the obsolete
# "FCF" column (U) for every year, and clears several stray figures that
# should not have been present for the estimate years (rows 7-9: D/Q/R/S/T/U/W/X).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 279272
$ws.Range("E2").Value = 12434
$ws.Range("F2").Value = 12434
$ws.Range("G2").Value = 12526
$ws.Range("H2").Value = 9798
$ws.Range("I2").Value = 9377
$ws.Range("J2").Value = 421
$ws.Range("K2").Value = 3155482
$ws.Range("L2").Value = 2936546
$ws.Range("M2").Value = 218936
$ws.Range("N2").Value = 211273
$ws.Range("O2").Value = 7663
$ws.Range("P2").Value = 14495
$ws.Range("Q2").Value = 50566
$ws.Range("R2").Value = -49853
$ws.Range("S2").Value = -4894
$ws.Range("T2").Value = 3176
$ws.Range("V2").Value = 461411
$ws.Range("W2").Value = 4.45
$ws.Range("X2").Value = 3.51
$ws.Range("Y2").Value = 4.57
$ws.Range("Z2").Value = 0.32
$ws.Range("AA2").Value = 1341.28
$ws.Range("AB2").Value = 1410.46
$ws.Range("AC2").Value = 3235
$ws.Range("AD2").Value = 9.890000000000001
$ws.Range("AE2").Value = 72879
$ws.Range("AF2").Value = 0.44
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 1.88
$ws.Range("AI2").Value = 18.55
$ws.Range("AJ2").Value = 289894062
$ws.Range("U2").ClearContents()

# Row 3
$ws.Range("D3").Value = 340237
$ws.Range("E3").Value = 10105
$ws.Range("F3").Value = 10105
$ws.Range("G3").Value = 11964
$ws.Range("H3").Value = 9543
$ws.Range("I3").Value = 9097
$ws.Range("J3").Value = 445
$ws.Range("K3").Value = 3269127
$ws.Range("L3").Value = 3039410
$ws.Range("M3").Value = 229717
$ws.Range("N3").Value = 221497
$ws.Range("O3").Value = 8221
$ws.Range("P3").Value = 14800
$ws.Range("Q3").Value = 25255
$ws.Range("R3").Value = -31681
$ws.Range("S3").Value = -8272
$ws.Range("T3").Value = 5035
$ws.Range("V3").Value = 453719
$ws.Range("W3").Value = 2.97
$ws.Range("X3").Value = 2.8
$ws.Range("Y3").Value = 4.2
$ws.Range("Z3").Value = 0.3
$ws.Range("AA3").Value = 1323.11
$ws.Range("AB3").Value = 1452.13
$ws.Range("AC3").Value = 3093
$ws.Range("AD3").Value = 7.63
$ws.Range("AE3").Value = 74829
$ws.Range("AF3").Value = 0.32
$ws.Range("AG3").Value = 650
$ws.Range("AH3").Value = 2.75
$ws.Range("AI3").Value = 21.15
$ws.Range("AJ3").Value = 296003062
$ws.Range("U3").ClearContents()

# Row 4
$ws.Range("D4").Value = 345998
$ws.Range("E4").Value = 16141
$ws.Range("F4").Value = 16141
$ws.Range("G4").Value = 18200
$ws.Range("H4").Value = 13997
$ws.Range("I4").Value = 13305
$ws.Range("J4").Value = 693
$ws.Range("K4").Value = 3481775
$ws.Range("L4").Value = 3247875
$ws.Range("M4").Value = 233900
$ws.Range("N4").Value = 224877
$ws.Range("O4").Value = 9023
$ws.Range("P4").Value = 14800
$ws.Range("Q4").Value = 34354
$ws.Range("R4").Value = -15753
$ws.Range("S4").Value = -13217
$ws.Range("T4").Value = 5577
$ws.Range("V4").Value = 447387
$ws.Range("W4").Value = 4.67
$ws.Range("X4").Value = 4.05
$ws.Range("Y4").Value = 5.96
$ws.Range("Z4").Value = 0.42
$ws.Range("AA4").Value = 1388.57
$ws.Range("AB4").Value = 1480.39
$ws.Range("AC4").Value = 4495
$ws.Range("AD4").Value = 6.95
$ws.Range("AE4").Value = 75971
$ws.Range("AF4").Value = 0.41
$ws.Range("AG4").Value = 1050
$ws.Range("AH4").Value = 3.36
$ws.Range("AI4").Value = 23.36
$ws.Range("AJ4").Value = 296003062
$ws.Range("U4").ClearContents()

# Row 5
$ws.Range("D5").Value = 391711
$ws.Range("E5").Value = 27181
$ws.Range("F5").Value = 27181
$ws.Range("G5").Value = 27987
$ws.Range("H5").Value = 21166
$ws.Range("I5").Value = 20368
$ws.Range("J5").Value = 798
$ws.Range("K5").Value = 3600893
$ws.Range("L5").Value = 3352614
$ws.Range("M5").Value = 248279
$ws.Range("N5").Value = 238737
$ws.Range("O5").Value = 9543
$ws.Range("P5").Value = 14800
$ws.Range("Q5").Value = 60119
$ws.Range("R5").Value = -102530
$ws.Range("S5").Value = 40918
$ws.Range("T5").Value = 7403
$ws.Range("V5").Value = 482696
$ws.Range("W5").Value = 6.94
$ws.Range("X5").Value = 5.4
$ws.Range("Y5").Value = 8.789999999999999
$ws.Range("Z5").Value = 0.6
$ws.Range("AA5").Value = 1350.34
$ws.Range("AB5").Value = 1577.55
$ws.Range("AC5").Value = 6881
$ws.Range("AD5").Value = 7.24
$ws.Range("AE5").Value = 80654
$ws.Range("AF5").Value = 0.62
$ws.Range("AG5").Value = 1550
$ws.Range("AH5").Value = 3.11
$ws.Range("AI5").Value = 22.53
$ws.Range("AJ5").Value = 296003062
$ws.Range("U5").ClearContents()

# Row 6
$ws.Range("D6").Value = 325159
$ws.Range("E6").Value = 31522
$ws.Range("F6").Value = 31522
$ws.Range("G6").Value = 31497
$ws.Range("H6").Value = 22752
$ws.Range("I6").Value = 22333
$ws.Range("K6").Value = 3850086
$ws.Range("L6").Value = 3579002
$ws.Range("M6").Value = 271085
$ws.Range("N6").Value = 264292
$ws.Range("P6").Value = 15012
$ws.Range("Q6").Value = -54317
$ws.Range("R6").Value = -16107
$ws.Range("S6").Value = 38672
$ws.Range("T6").Value = 4971
$ws.Range("V6").Value = 554632
$ws.Range("W6").Value = 9.69
$ws.Range("X6").Value = 7
$ws.Range("Y6").Value = 8.880000000000001
$ws.Range("Z6").Value = 0.61
$ws.Range("AA6").Value = 1320.25
$ws.Range("AB6").Value = 1705.77
$ws.Range("AC6").Value = 7457
$ws.Range("AD6").Value = 4.86
$ws.Range("AE6").Value = 88026
$ws.Range("AF6").Value = 0.41
$ws.Range("AG6").Value = 1900
$ws.Range("AH6").Value = 5.24
$ws.Range("AI6").Value = 25.54
$ws.Range("AJ6").Value = 300242062
$ws.Range("U6").ClearContents()

# Row 7
$ws.Range("E7").Value = 32423
$ws.Range("G7").Value = 34242
$ws.Range("H7").Value = 24885
$ws.Range("I7").Value = 24791
$ws.Range("K7").Value = 4168198
$ws.Range("L7").Value = 3874324
$ws.Range("M7").Value = 293521
$ws.Range("N7").Value = 285614
$ws.Range("P7").Value = 14929
$ws.Range("Y7").Value = 9.02
$ws.Range("Z7").Value = 0.62
$ws.Range("AA7").Value = 1319.95
$ws.Range("AC7").Value = 8257
$ws.Range("AD7").Value = 4.02
$ws.Range("AE7").Value = 97960
$ws.Range("AF7").Value = 0.34
$ws.Range("AG7").Value = 2135
$ws.Range("AH7").Value = 6.43
$ws.Range("AI7").Value = 25.85
$ws.Range("D7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()

# Row 8
$ws.Range("E8").Value = 33647
$ws.Range("G8").Value = 32527
$ws.Range("H8").Value = 24028
$ws.Range("I8").Value = 23661
$ws.Range("K8").Value = 4335910
$ws.Range("L8").Value = 4024220
$ws.Range("M8").Value = 311402
$ws.Range("N8").Value = 303706
$ws.Range("P8").Value = 14929
$ws.Range("Y8").Value = 8.029999999999999
$ws.Range("Z8").Value = 0.5600000000000001
$ws.Range("AA8").Value = 1292.29
$ws.Range("AC8").Value = 7881
$ws.Range("AD8").Value = 4.21
$ws.Range("AE8").Value = 104165
$ws.Range("AF8").Value = 0.32
$ws.Range("AG8").Value = 2177
$ws.Range("AH8").Value = 6.56
$ws.Range("AI8").Value = 27.63
$ws.Range("D8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()

# Row 9
$ws.Range("E9").Value = 35726
$ws.Range("G9").Value = 34394
$ws.Range("H9").Value = 25282
$ws.Range("I9").Value = 24731
$ws.Range("K9").Value = 4491275
$ws.Range("L9").Value = 4158932
$ws.Range("M9").Value = 332343
$ws.Range("N9").Value = 322665
$ws.Range("P9").Value = 14980
$ws.Range("Y9").Value = 7.9
$ws.Range("Z9").Value = 0.57
$ws.Range("AA9").Value = 1251.4
$ws.Range("AC9").Value = 8237
$ws.Range("AD9").Value = 4.03
$ws.Range("AE9").Value = 110667
$ws.Range("AF9").Value = 0.3
$ws.Range("AG9").Value = 2317
$ws.Range("AH9").Value = 6.98
$ws.Range("AI9").Value = 28.13
$ws.Range("D9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
